$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "TOSHIBA 43'' Smart LED TV, ..." -> split the run into
#   "TOSHIBA 43" + " Smart LED TV, ..." (the '' is removed) as two <w:r>
#   elements sharing identical run formatting.
# ---------------------------------------------------------------------------

# First, drop the stray '' right after "TOSHIBA 43" (inches mark).
$r1 = $d.Content
$found1 = $r1.Find.Execute("TOSHIBA 43''", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "TOSHIBA 43", 2)

# Now force a run split exactly at the new word boundary by toggling a
# character-formatting property on just that sub-range and reverting it.
# Applying (then clearing) the property is what makes Word materialize a
# fresh <w:r> at this position instead of silently merging back into its
# neighbour.
$r2 = $d.Content
$found2 = $r2.Find.Execute("TOSHIBA 43", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
$r2.Font.Bold = 1
$r2.Font.Bold = 0

# ---------------------------------------------------------------------------
# Edit 2: move the "_GoBack" bookmark so it wraps the
#   "Toshiba LED 32 Inch HD TV with Built-In Receiver" run instead of sitting
#   as an empty bookmark after it.
# ---------------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Toshiba LED 32 Inch HD TV with Built-In Receiver*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Leelawadee" w:hAnsi="Leelawadee" w:cs="Leelawadee"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:rPr><w:rFonts w:ascii="Leelawadee" w:hAnsi="Leelawadee" w:cs="Leelawadee"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>Toshiba LED 32 Inch HD TV with Built-In Receiver</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    $target.Range.InsertXML($xml)
}

Write-Host "Edits applied."
